# "Edited until Nehemiah 2:10" — mark the newly finished books as done.
#
# The "all" sheet tracks Bible-reading progress: column F ("Book done") is a
# manual flag per book; column E ("Verses done") is a formula
# (=IF(F=1,C,0)) that only counts a book's verses once it's flagged done.
# Flip F to 1 for every book that has now been completed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# Row -> Book just finished
#  16 -> Ezra
#  18 -> Esther
#  27 -> Ezekiel
#  28 -> Daniel
#  37 -> Zephaniah
#  38 -> Haggai
#  39 -> Zechariah
#  40 -> Malachi
$doneRows = @(16, 18, 27, 28, 37, 38, 39, 40)
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 6).Value = 1
}

# Leave the selection where work left off.
$ws.Range("A2:E40").Select()
